$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Populate the login table
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# View adjustments
$excel.ActiveWindow.Zoom = 235
$ws.Range("C2").Select()
